$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new data row (row 54) with the same look & feel as the previous rows.

# A54 holds a date-like string ("07-11-2025") that must stay literal text,
# not get auto-converted into a date serial number. Build it as a formula
# that returns a text string, then convert the formula to its static value
# in place - this keeps the cell type as a plain text/shared-string cell.
$ws.Range("A54").Formula = "=""07-11-2025"""
$ws.Range("A54").Copy()
$ws.Range("A54").PasteSpecial(-4163)
$excel.CutCopyMode = 0

# B54 is an ordinary text value.
$ws.Range("B54").Value = "The price of gold in India today is ₹12,202 per gram for 24 karat gold, ₹11,185 per gram for 22 karat gold and ₹9,152 per gram for 18 karat gold (also called 999 gold)."

# Copy the cell formatting (borders/fill/wrap) from the row above so the new
# row matches the rest of the table's style.
$ws.Range("A53:B53").Copy()
$ws.Range("A54:B54").PasteSpecial(-4122)
$excel.CutCopyMode = 0
